$d = $word.ActiveDocument

$newText = "Ημερομηνίες παρατήρησης για τον αστερισμό του Pegasus: 8-17 Οκτωβρίου, 7-16 Νοεμβρίου,"

$updated = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "2018*?ερσε?*") {
        $r = $p.Range
        # Exclude the trailing paragraph-mark character from the range.
        [void]$r.MoveEnd(1, -1)
        # Delete all existing runs (and their formatting) in this range ...
        $r.Text = ""
        # ... then insert the replacement as a brand-new, unformatted run,
        # matching the target (no <w:rPr>) rather than inheriting the
        # formatting of any of the deleted runs.
        [void]$r.InsertAfter($newText)
        $updated = $updated + 1
    }
}

Write-Host "Paragraphs updated:" $updated
